$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5333.3335
$ws.Range("I32").Value = 5499.5
$ws.Range("J32").Value = 5001
$ws.Range("K32").Value = 5499.5
$ws.Range("L32").Value = 5001
$ws.Range("M32").Value = -5173.5
$ws.Range("N32").Value = -5653
$ws.Range("H62").Value = 22733980
$ws.Range("I62").Value = 35720384
$ws.Range("K62").Value = 35720384
$ws.Range("M62").Value = -35719760
$ws.Range("H65").Value = 22733980
$ws.Range("I65").Value = 35720384
$ws.Range("K65").Value = 178601920
$ws.Range("M65").Value = -178598800
$ws.Range("H98").Value = 3365099.5
$ws.Range("I98").Value = 3368585.2
$ws.Range("J98").Value = 3333726.2
$ws.Range("K98").Value = 3368585.2
$ws.Range("L98").Value = 3333726.2
$ws.Range("M98").Value = -3367087.2
$ws.Range("N98").Value = -3336722.2
$ws.Range("H106").Value = 1748.5
$ws.Range("I106").Value = 1712.5714
$ws.Range("K106").Value = 1712.5714
$ws.Range("M106").Value = -1081.5714
$ws.Range("H107").Value = 50004900
$ws.Range("I107").Value = 31255400
$ws.Range("J107").Value = 125002900
$ws.Range("K107").Value = 31255400
$ws.Range("L107").Value = 125002900
$ws.Range("M107").Value = -31253480
$ws.Range("N107").Value = -125006740
$ws.Range("H122").Value = 3365099.5
$ws.Range("I122").Value = 3368585.2
$ws.Range("J122").Value = 3333726.2
$ws.Range("K122").Value = 10105755.6
$ws.Range("L122").Value = 10001178.6
$ws.Range("M122").Value = -10103305.6
$ws.Range("N122").Value = -10006078.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 3701.2
$ws.Range("J29").Value = 3502
$ws.Range("L29").Value = 3502
$ws.Range("N29").Value = -4118
$ws.Range("H32").Value = 6866.9214
$ws.Range("I32").Value = 3393.7046
$ws.Range("K32").Value = 3393.7046
$ws.Range("M32").Value = -3106.7046
$ws.Range("H102").Value = 2850.3333
$ws.Range("I102").Value = 2655
$ws.Range("K102").Value = 2655
$ws.Range("M102").Value = -1033

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 52496.332
$ws.Range("J74").Value = 52496.332
$ws.Range("L74").Value = 52496.332
$ws.Range("N74").Value = -54368.332
$ws.Range("H77").Value = 52496.332
$ws.Range("J77").Value = 52496.332
$ws.Range("L77").Value = 157488.996
$ws.Range("N77").Value = -166848.996
$ws.Range("H86").Value = 2950.8928
$ws.Range("I86").Value = 3567.0715
$ws.Range("J86").Value = 2334.7144
$ws.Range("K86").Value = 3567.0715
$ws.Range("L86").Value = 2334.7144
$ws.Range("M86").Value = -2444.0715
$ws.Range("N86").Value = -4580.7144
$ws.Range("H89").Value = 2950.8928
$ws.Range("I89").Value = 3567.0715
$ws.Range("J89").Value = 2334.7144
$ws.Range("K89").Value = 17835.3575
$ws.Range("L89").Value = 11673.572
$ws.Range("M89").Value = -12219.3575
$ws.Range("N89").Value = -22905.572
$ws.Range("H105").Value = 2058.25
$ws.Range("I105").Value = 1931.5
$ws.Range("J105").Value = 3199
$ws.Range("K105").Value = 1931.5
$ws.Range("L105").Value = 3199
$ws.Range("M105").Value = -184.5
$ws.Range("N105").Value = -6693
$ws.Range("H107").Value = 3983
$ws.Range("I107").Value = 3983
$ws.Range("K107").Value = 3983
$ws.Range("M107").Value = -2063

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 10640.5
$ws.Range("J28").Value = 10640.5
$ws.Range("L28").Value = 10640.5
$ws.Range("N28").Value = -11130.5
$ws.Range("H58").Value = 5085.517
$ws.Range("I58").Value = 4758.9546
$ws.Range("J58").Value = 6111.857
$ws.Range("K58").Value = 4758.9546
$ws.Range("L58").Value = 6111.857
$ws.Range("M58").Value = -4555.9546
$ws.Range("N58").Value = -6517.857
$ws.Range("H96").Value = 11783
$ws.Range("J96").Value = 11783
$ws.Range("L96").Value = 11783
$ws.Range("N96").Value = -17275
$ws.Range("H106").Value = 20555.5
$ws.Range("J106").Value = 20555.5
$ws.Range("L106").Value = 20555.5
$ws.Range("N106").Value = -23079.5
$ws.Range("H107").Value = 1200.2084
$ws.Range("I107").Value = 733
$ws.Range("J107").Value = 1978.8889
$ws.Range("K107").Value = 733
$ws.Range("L107").Value = 1978.8889
$ws.Range("M107").Value = 1187
$ws.Range("N107").Value = -5818.8889
$ws.Range("H132").Value = 2447.4736
$ws.Range("I132").Value = 2499.077
$ws.Range("K132").Value = 7497.231000000001
$ws.Range("M132").Value = -4967.231000000001
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060
$ws.Range("H134").Value = 9154.8125
$ws.Range("I134").Value = 9431.799999999999
$ws.Range("K134").Value = 28295.4
$ws.Range("M134").Value = -25760.4
$ws.Range("H136").Value = 5085.517
$ws.Range("I136").Value = 4758.9546
$ws.Range("J136").Value = 6111.857
$ws.Range("K136").Value = 14276.8638
$ws.Range("L136").Value = 18335.571
$ws.Range("M136").Value = -11726.8638
$ws.Range("N136").Value = -23435.571

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 252.75
$ws.Range("I6").Value = 79.5
$ws.Range("K6").Value = 238.5
$ws.Range("M6").Value = -125.5
$ws.Range("H26").Value = 4083.75
$ws.Range("I26").Value = 390.33334
$ws.Range("K26").Value = 1171.00002
$ws.Range("M26").Value = -883.0000199999999
$ws.Range("H99").Value = 8425.0625
$ws.Range("I99").Value = 3724.25
$ws.Range("K99").Value = 11172.75
$ws.Range("M99").Value = -8926.75
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H113").Value = 646.7857
$ws.Range("I113").Value = 246.5
$ws.Range("J113").Value = 713.5
$ws.Range("K113").Value = 739.5
$ws.Range("L113").Value = 2140.5
$ws.Range("M113").Value = 1430.5
$ws.Range("N113").Value = -6480.5
$ws.Range("H132").Value = 6275.36
$ws.Range("I132").Value = 6429.7827
$ws.Range("K132").Value = 57868.04429999999
$ws.Range("M132").Value = -55338.04429999999
$ws.Range("H140").Value = 3086.9412
$ws.Range("I140").Value = 1843.375
$ws.Range("J140").Value = 4192.3335
$ws.Range("K140").Value = 5530.125
$ws.Range("L140").Value = 12577.0005
$ws.Range("M140").Value = -350.125
$ws.Range("N140").Value = -22937.0005

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4204.826
$ws.Range("I70").Value = 4245.6733
$ws.Range("J70").Value = 4104.75
$ws.Range("K70").Value = 4245.6733
$ws.Range("L70").Value = 4104.75
$ws.Range("M70").Value = -3975.6733
$ws.Range("N70").Value = -4644.75
$ws.Range("H73").Value = 4204.826
$ws.Range("I73").Value = 4245.6733
$ws.Range("J73").Value = 4104.75
$ws.Range("K73").Value = 4245.6733
$ws.Range("L73").Value = 4104.75
$ws.Range("M73").Value = -3309.6733
$ws.Range("N73").Value = -5976.75
$ws.Range("H99").Value = 17566.75
$ws.Range("I99").Value = 17566.75
$ws.Range("K99").Value = 17566.75
$ws.Range("M99").Value = -15320.75
$ws.Range("H132").Value = 26471.637
$ws.Range("I132").Value = 28675.4
$ws.Range("J132").Value = 4434
$ws.Range("K132").Value = 86026.20000000001
$ws.Range("L132").Value = 13302
$ws.Range("M132").Value = -83496.20000000001
$ws.Range("N132").Value = -18362

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7724.722
$ws.Range("I7").Value = 7724.722
$ws.Range("K7").Value = 7724.722
$ws.Range("M7").Value = -7612.722
$ws.Range("H20").Value = 8437.733
$ws.Range("I20").Value = 6109.8335
$ws.Range("K20").Value = 6109.8335
$ws.Range("M20").Value = -5883.8335
$ws.Range("H64").Value = 23075
$ws.Range("J64").Value = 23075
$ws.Range("L64").Value = 23075
$ws.Range("N64").Value = -23525
$ws.Range("H67").Value = 23075
$ws.Range("J67").Value = 23075
$ws.Range("L67").Value = 23075
$ws.Range("N67").Value = -24635
$ws.Range("H122").Value = 4205.125
$ws.Range("I122").Value = 3461.5
$ws.Range("K122").Value = 10384.5
$ws.Range("M122").Value = -7934.5
$ws.Range("H126").Value = 7724.722
$ws.Range("I126").Value = 7724.722
$ws.Range("K126").Value = 23174.166
$ws.Range("M126").Value = -20704.166

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 19714
$ws.Range("I14").Value = 16999.666
$ws.Range("J14").Value = 36000
$ws.Range("K14").Value = 16999.666
$ws.Range("L14").Value = 36000
$ws.Range("M14").Value = -16831.666
$ws.Range("N14").Value = -36336
$ws.Range("H62").Value = 9763.556
$ws.Range("I62").Value = 6999.5
$ws.Range("J62").Value = 10553.286
$ws.Range("K62").Value = 6999.5
$ws.Range("L62").Value = 10553.286
$ws.Range("M62").Value = -6375.5
$ws.Range("N62").Value = -11801.286
$ws.Range("H63").Value = 19550
$ws.Range("I63").Value = 10000
$ws.Range("K63").Value = 10000
$ws.Range("M63").Value = -9376
$ws.Range("H65").Value = 9763.556
$ws.Range("I65").Value = 6999.5
$ws.Range("J65").Value = 10553.286
$ws.Range("K65").Value = 34997.5
$ws.Range("L65").Value = 52766.43
$ws.Range("M65").Value = -31877.5
$ws.Range("N65").Value = -59006.43
$ws.Range("H66").Value = 19550
$ws.Range("I66").Value = 10000
$ws.Range("K66").Value = 30000
$ws.Range("M66").Value = -26880
$ws.Range("H122").Value = 13319.259
$ws.Range("I122").Value = 3097.348
$ws.Range("J122").Value = 72095.25
$ws.Range("K122").Value = 9292.044
$ws.Range("L122").Value = 216285.75
$ws.Range("M122").Value = -6842.044
$ws.Range("N122").Value = -221185.75
$ws.Range("H132").Value = 1785.5
$ws.Range("I132").Value = 1569.3334
$ws.Range("J132").Value = 2001.6666
$ws.Range("K132").Value = 4708.0002
$ws.Range("L132").Value = 6004.9998
$ws.Range("M132").Value = -2178.0002
$ws.Range("N132").Value = -11064.9998
$ws.Range("H136").Value = 436471.1
$ws.Range("I136").Value = 455628.88
$ws.Range("K136").Value = 1366886.64
$ws.Range("M136").Value = -1364336.64
